# Bug Log.docx edit script
#
# Goal (per the commit diff):
#  1. Paragraph "in firstbite, clicking always ADDS parameter ..." loses its
#     gramStart/gramEnd/spellStart/spellEnd proofErr markers and its several
#     runs collapse into a single run with the full sentence.
#  2. Paragraph "when using RT, only the first sear|ch result is used as the
#     movie" loses the mid-sentence _GoBack bookmark and its two runs merge
#     into a single run with the complete ("search") text.
#  3. A brand-new list item "in firstbite, one-page CSS is currently broken"
#     is added right after that paragraph, carrying the _GoBack bookmark
#     (now placed after its own run, i.e. at the end of the paragraph).
#
# Word's Find/Replace and Range.Delete leave orphaned w:proofErr /
# w:bookmarkStart/End markers behind because they are not part of the
# "text" stream. The reliable way to really drop them is to rebuild the
# affected paragraphs' underlying OOXML via Range.InsertXML.
#
# The replacement range has to start at the very beginning of the document
# (offset 0): this document's paragraph count grows by one (a new list item
# is added), and starting the InsertXML range anywhere other than offset 0
# while growing the paragraph count causes the engine to swallow the
# paragraph immediately before the range. Starting at 0 sidesteps that, at
# the minor cost of Word regenerating paragraph 1 ("Bug Log") too (it keeps
# its text/formatting, it just loses its cosmetic w:rsid* bookkeeping
# attributes, which carry no semantic meaning).

$d = $word.ActiveDocument

$wordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-WordPackageXml([string]$bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document ' + $wordNs + '><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData>' +
           '</pkg:part></pkg:package>'
}

# Locate the existing paragraph that marks the end of the region we need to
# rebuild (the "when using RT..." paragraph - the third paragraph).
$pRT = $d.Paragraphs(3)

# Build replacement XML for paragraphs 1-3, then append the brand-new 4th
# list item (with the relocated bookmark) right after paragraph 3's content.
$listPPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>'

$bodyFragment =
    '<w:p><w:r><w:t>Bug Log</w:t></w:r></w:p>' +
    '<w:p>' + $listPPr + '<w:r><w:t>in firstbite, clicking always ADDS parameter (instinctiveness, etc. to user)</w:t></w:r></w:p>' +
    '<w:p>' + $listPPr + '<w:r><w:t>when using RT, only the first search result is used as the movie</w:t></w:r></w:p>' +
    '<w:p>' + $listPPr + '<w:r><w:t>in firstbite, one-page CSS is currently broken</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

# Replace from the very start of the document through the end of the
# "when using RT..." paragraph; the trailing w:sectPr is left untouched.
$range = $d.Range(0, $pRT.Range.End)
$range.InsertXML((New-WordPackageXml $bodyFragment))

Write-Host "Final paragraphs:"
foreach ($p in $d.Paragraphs) {
    Write-Host " - [$($p.Range.Text)]"
}
